$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '25.761.22'
$ws.Range('E2').Value = '  -0.99%  '
# Row 3
$ws.Range('D3').Value = '1.627.15'
$ws.Range('E3').Value = '  -0.92%  '
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.43'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.41%  '
# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5108'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.41%  '
# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.03%  '
# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2562'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.01%  '
# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06321'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.41%  '
# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.44'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.42%  '
# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07780'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.25%  '
# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.656.43'
$ws.Range('E12').Value = '  +0.87%  '
# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.233'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.25%  '
# Row 14
$ws.Range('D14').Value = '1.850.30'
$ws.Range('E14').Value = '  -1.06%  '
# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5517'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.68%  '
# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.64'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.87%  '
# Row 17
$ws.Range('D17').Value = '0.0₅7548'
$ws.Range('E17').Value = '  -2.00%  '
# Row 18
$ws.Range('D18').Value = '25.783.39'
$ws.Range('E18').Value = '  -0.99%  '
# Row 19
$ws.Range('E19').Value = '  +0.06%  '
# Row 20
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '194.11'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.16%  '
# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.405'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.41%  '
# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.837'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.76%  '
# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.000'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.57%  '
# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.16%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.886'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.87%  '
# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '141.98'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.50%  '
# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1257'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.62%  '
# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.55'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.29%  '
# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.741'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.84%  '
# Row 30
$ws.Range('E30').Value = '  +0.31%  '
# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.04881'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.38%  '
# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.227'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.72%  '
# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.170'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.34%  '
# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.537'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.89%  '
# Row 35
$ws.Range('E35').Value = '  +0.36%  '
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.8930'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.61%  '
# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.5504'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.04%  '
# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.531'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.00%  '
# Row 39
$ws.Range('D39').Value = '1.111.54'
$ws.Range('E39').Value = '  -2.51%  '
# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01549'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.75%  '
# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.000'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.16%  '
# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.553'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.14%  '
# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.7967'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.72%  '
# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '97.15'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.12%  '
# Row 45
$ws.Range('D45').Value = '1.773.76'
$ws.Range('E45').Value = '  -0.48%  '
# Row 46
$ws.Range('E46').Value = '  -13.59%  '
# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4429'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.17%  '
# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.002'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.27%  '
# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '54.59'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.57%  '
# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05134'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.55%  '
# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.551'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.07%  '
